$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row above the old row 2 (the thin separator row). This
#    shifts the old rows 2-5 down to become rows 3-6.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).Insert()

# ---------------------------------------------------------------------------
# 2. Fill the new row 2 with the unit-of-measure captions
#    (Kyrgyz / Russian / English), italic 9pt Times New Roman, centered.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,1).Value = "(бирдик)"
$ws.Cells.Item(2,2).Value = "(единиц)"
$ws.Cells.Item(2,3).Value = "(units)"

$newRow2 = $ws.Range("A2:C2")
$newRow2.WrapText = $false
$newRow2.Font.Size = 9
$newRow2.Font.Bold = $false
$newRow2.Font.Italic = $true
$newRow2.HorizontalAlignment = -4108
$newRow2.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. Row heights / column widths
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 41.25
$ws.Rows.Item(6).RowHeight = 18
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 35

# ---------------------------------------------------------------------------
# 4. Extend the thin separator row (old row 2, now row 3) across J:K by
#    cloning the formatting already present in column I of that row.
# ---------------------------------------------------------------------------
$ws.Range("I3").Copy()
$ws.Range("J3:K3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. Add the 2022 / 2023 year headers (row 4) - clone column I formatting.
# ---------------------------------------------------------------------------
$ws.Range("I4").Copy()
$ws.Range("J4:K4").PasteSpecial(-4122)
$ws.Cells.Item(4,10).Value = 2022
$ws.Cells.Item(4,11).Value = 2023

# ---------------------------------------------------------------------------
# 6. Row 5 ("Voluntarily surrendered firearms") - add J5 ("-") and K5 (219)
# ---------------------------------------------------------------------------
$ws.Range("I5").Copy()
$ws.Range("J5:K5").PasteSpecial(-4122)
$ws.Cells.Item(5,10).Value = "-"
$ws.Cells.Item(5,10).HorizontalAlignment = -4152
$ws.Cells.Item(5,11).Value = 219

# ---------------------------------------------------------------------------
# 7. Row 6 ("Seized firearms") - add J6 ("-") and K6 (171)
# ---------------------------------------------------------------------------
$ws.Range("I6").Copy()
$ws.Range("J6:K6").PasteSpecial(-4122)
$ws.Cells.Item(6,10).Value = "-"
$ws.Cells.Item(6,10).HorizontalAlignment = -4152
$ws.Cells.Item(6,11).Value = 171

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 8. Page setup - print resolution
# ---------------------------------------------------------------------------
$ws.PageSetup.PrintQuality = 300

$ws.Cells.Item(1,1).Select()
